$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF") with the same formatting as the
# existing header cells (bold, bordered, centered - style used by B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I and J columns (rows 2-65) with their data values.
$newData = @(
    @{ Row = 2; I = 9; J = 9 },
    @{ Row = 3; I = 9; J = 9 },
    @{ Row = 4; I = 9; J = 9 },
    @{ Row = 5; I = 9; J = 9 },
    @{ Row = 6; I = 9; J = 9 },
    @{ Row = 7; I = 10; J = 10 },
    @{ Row = 8; I = 10; J = 10 },
    @{ Row = 9; I = 9; J = 9 },
    @{ Row = 10; I = 9; J = 9 },
    @{ Row = 11; I = 9; J = 9 },
    @{ Row = 12; I = 9; J = 9 },
    @{ Row = 13; I = 9; J = 9 },
    @{ Row = 14; I = 9; J = 9 },
    @{ Row = 15; I = 9; J = 9 },
    @{ Row = 16; I = 9; J = 9 },
    @{ Row = 17; I = 9; J = 9 },
    @{ Row = 18; I = 9; J = 9 },
    @{ Row = 19; I = 9; J = 9 },
    @{ Row = 20; I = 9; J = 9 },
    @{ Row = 21; I = 9; J = 9 },
    @{ Row = 22; I = 9; J = 9 },
    @{ Row = 23; I = 9; J = 9 },
    @{ Row = 24; I = 9; J = 9 },
    @{ Row = 25; I = 9; J = 9 },
    @{ Row = 26; I = 9; J = 9 },
    @{ Row = 27; I = 9; J = 9 },
    @{ Row = 28; I = 9; J = 9 },
    @{ Row = 29; I = 9; J = 9 },
    @{ Row = 30; I = 9; J = 9 },
    @{ Row = 31; I = 9; J = 9 },
    @{ Row = 32; I = 9; J = 9 },
    @{ Row = 33; I = 9; J = 9 },
    @{ Row = 34; I = 8; J = 8 },
    @{ Row = 35; I = 9; J = 9 },
    @{ Row = 36; I = 9; J = 9 },
    @{ Row = 37; I = 9; J = 9 },
    @{ Row = 38; I = 9; J = 9 },
    @{ Row = 39; I = 9; J = 9 },
    @{ Row = 40; I = 8; J = 9 },
    @{ Row = 41; I = 9; J = 9 },
    @{ Row = 42; I = 9; J = 9 },
    @{ Row = 43; I = 9; J = 9 },
    @{ Row = 44; I = 10; J = 10 },
    @{ Row = 45; I = 8; J = 9 },
    @{ Row = 46; I = 9; J = 9 },
    @{ Row = 47; I = 10; J = 10 },
    @{ Row = 48; I = 9; J = 9 },
    @{ Row = 49; I = 9; J = 9 },
    @{ Row = 50; I = 9; J = 9 },
    @{ Row = 51; I = 9; J = 9 },
    @{ Row = 52; I = 9; J = 9 },
    @{ Row = 53; I = 9; J = 10 },
    @{ Row = 54; I = 8; J = 8 },
    @{ Row = 55; I = 9; J = 9 },
    @{ Row = 56; I = 9; J = 9 },
    @{ Row = 57; I = 9; J = 9 },
    @{ Row = 58; I = 9; J = 9 },
    @{ Row = 59; I = 8; J = 9 },
    @{ Row = 60; I = 5; J = 5 },
    @{ Row = 61; I = 6; J = 6 },
    @{ Row = 62; I = 5; J = 5 },
    @{ Row = 63; I = 5; J = 5 },
    @{ Row = 64; I = 3; J = 3 },
    @{ Row = 65; I = 3; J = 3 }
)

foreach ($item in $newData) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}

# Update the used-range dimension to reflect the new columns.
$ws.UsedRange | Out-Null
